# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
# These two sheets list the same events, so the same cells are updated on both.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F2").Value = 13789
        $ws.Range("F6").Value = 498
        $ws.Range("F8").Value = 1030
        $ws.Range("F10").Value = 14749
        $ws.Range("F29").Value = 5412
        $ws.Range("F32").Value = 259
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F2").Value = 13789
        $ws.Range("F7").Value = 498
        $ws.Range("F9").Value = 1030
        $ws.Range("F11").Value = 14749
        $ws.Range("F30").Value = 5412
        $ws.Range("F33").Value = 259
    }
}
